$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '50.914.60'
$ws.Range("E2").Value = '  -0.49%  '

# Row 3
$ws.Range("D3").Value = '2.935.76'
$ws.Range("E3").Value = '  -0.77%  '

# Row 4
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("D5").Value = '''378.02'
$ws.Range("E5").Value = '  -0.35%  '

# Row 6
$ws.Range("D6").Value = '''100.63'
$ws.Range("E6").Value = '  -1.85%  '

# Row 7
$ws.Range("D7").Value = '''0.538'
$ws.Range("E7").Value = '  -0.02%  '

# Row 8
$ws.Range("E8").Value = '  +0.02%  '

# Row 9
$ws.Range("D9").Value = '''0.579'
$ws.Range("E9").Value = '  -1.52%  '

# Row 10
$ws.Range("D10").Value = '''35.94'
$ws.Range("E10").Value = '  -1.60%  '

# Row 11
$ws.Range("E11").Value = '  -0.45%  '

# Row 12
$ws.Range("D12").Value = '''0.0851'
$ws.Range("E12").Value = '  +1.44%  '

# Row 13
$ws.Range("D13").Value = '3.398.35'
$ws.Range("E13").Value = '  -0.65%  '

# Row 14
$ws.Range("B14").Value = 'Uniswap'
$ws.Range("C14").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D14").Value = '''12.23'
$ws.Range("E14").Value = '  +70.60%  '

# Row 15
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '''18.19'
$ws.Range("E15").Value = '  +1.28%  '

# Row 16
$ws.Range("D16").Value = '''7.61'
$ws.Range("E16").Value = '  +3.48%  '

# Row 17
$ws.Range("D17").Value = '2.933.11'
$ws.Range("E17").Value = '  -0.43%  '

# Row 18
$ws.Range("D18").Value = '''0.993'
$ws.Range("E18").Value = '  +1.15%  '

# Row 19
$ws.Range("D19").Value = '50.896.30'
$ws.Range("E19").Value = '  -0.54%  '

# Row 20
$ws.Range("D20").Value = '''3.01'
$ws.Range("E20").Value = '  -6.48%  '

# Row 21
$ws.Range("D21").Value = '''12.38'
$ws.Range("E21").Value = '  -1.33%  '

# Row 22
$ws.Range("E22").Value = '  -0.56%  '

# Row 23
$ws.Range("D23").Value = '''69.46'
$ws.Range("E23").Value = '  +1.71%  '

# Row 24
$ws.Range("D24").Value = '''266.66'
$ws.Range("E24").Value = '  +2.02%  '

# Row 25
$ws.Range("D25").Value = '''3.24'
$ws.Range("E25").Value = '  +13.61%  '

# Row 26
$ws.Range("D26").Value = '''7.94'
$ws.Range("E26").Value = '  -3.46%  '

# Row 27
$ws.Range("E27").Value = '  -0.09%  '

# Row 28
$ws.Range("D28").Value = '''7.10'
$ws.Range("E28").Value = '  -6.56%  '

# Row 29
$ws.Range("D29").Value = '''25.55'
$ws.Range("E29").Value = '  -0.49%  '

# Row 30
$ws.Range("D30").Value = '''0.163'
$ws.Range("E30").Value = '  -2.91%  '

# Row 31
$ws.Range("E31").Value = '  -2.50%  '

# Row 32
$ws.Range("D32").Value = '''10.00'
$ws.Range("E32").Value = '  +2.25%  '

# Row 33
$ws.Range("D33").Value = '''50.51'
$ws.Range("E33").Value = '  -0.02%  '

# Row 34
$ws.Range("E34").Value = '  +0.25%  '

# Row 35
$ws.Range("D35").Value = '''33.46'
$ws.Range("E35").Value = '  -0.55%  '

# Row 36
$ws.Range("D36").Value = '''0.0432'
$ws.Range("E36").Value = '  -2.41%  '

# Row 37
$ws.Range("E37").Value = '  -0.04%  '

# Row 38
$ws.Range("E38").Value = '  +3.75%  '

# Row 39
$ws.Range("D39").Value = '''0.116'
$ws.Range("E39").Value = '  +0.97%  '

# Row 40
$ws.Range("D40").Value = '''16.52'
$ws.Range("E40").Value = '  -1.98%  '

# Row 41
$ws.Range("D41").Value = '''1.81'
$ws.Range("E41").Value = '  +1.90%  '

# Row 42
$ws.Range("D42").Value = '''2.48'
$ws.Range("E42").Value = '  -2.67%  '

# Row 43
$ws.Range("D43").Value = '''120.02'
$ws.Range("E43").Value = '  -1.53%  '

# Row 44
$ws.Range("D44").Value = '''21.16'

# Row 45
$ws.Range("D45").Value = '''3.43'
$ws.Range("E45").Value = '  +7.28%  '

# Row 46
$ws.Range("E46").Value = '  -1.71%  '

# Row 47
$ws.Range("E47").Value = '  -1.12%  '

# Row 48
$ws.Range("D48").Value = '2.008.37'
$ws.Range("E48").Value = '  +0.15%  '

# Row 49
$ws.Range("E49").Value = '  -4.90%  '

# Row 50
$ws.Range("D50").Value = '''0.0315'
$ws.Range("E50").Value = '  -5.41%  '

# Row 51
$ws.Range("D51").Value = '''5.28'
$ws.Range("E51").Value = '  +4.65%  '
